$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = "What the maximum number of headers I can display in my log?"
$ws.Range("B12").Value = "llama3.2:latest"
$ws.Range("C12").Value = "The maximum number of headers you can display in your log is up to 50."
